# Actualizacion Datos Personales 4 nov
# Updates statistics (Aprobados/Reprobados/Porcentajes/Promedio/Blancos) for
# a handful of groups on the "1er Parcial" and "3er Parcial" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 (e.g. 3AEV group)
    $ws.Range("E2").Value = 14
    $ws.Range("F2").Value = 19
    $ws.Range("G2").Value = 42.42
    $ws.Range("H2").Value = 57.58
    $ws.Range("I2").Value = 6.2
    $ws.Range("J2").Value = 0
    $ws.Range("K2").Value = 0

    # Row 3 (e.g. 3ASV group)
    $ws.Range("E3").Value = 13
    $ws.Range("F3").Value = 8
    $ws.Range("G3").Value = 61.9
    $ws.Range("H3").Value = 38.1
    $ws.Range("I3").Value = 7.2
    $ws.Range("J3").Value = 0
    $ws.Range("K3").Value = 0

    # Row 18 (e.g. 3BEM group)
    $ws.Range("E18").Value = 18
    $ws.Range("F18").Value = 13
    $ws.Range("G18").Value = 58.06
    $ws.Range("H18").Value = 41.94
    $ws.Range("I18").Value = 7.2
    $ws.Range("J18").Value = 13
    $ws.Range("K18").Value = 41.94

    # Row 23
    $ws.Range("E23").Value = 31
    $ws.Range("F23").Value = 5
    $ws.Range("G23").Value = 86.11
    $ws.Range("H23").Value = 13.89
    $ws.Range("I23").Value = 8.199999999999999
    $ws.Range("J23").Value = 5
    $ws.Range("K23").Value = 13.89
}
